$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.572.12'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.974.09'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.628'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.11'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.97%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0789'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.843'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').Value = '2.262.94'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.38'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').Value = '1.971.64'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '36.551.67'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.83'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = '0.0₃0856'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.37'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.146'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +7.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +19.29%  '
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0616'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.51'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.28'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.03%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('E39').Value = '  -13.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0973'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.61%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').Value = '1.365.03'
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.19%  '
$ws.Range('D51').Value = '2.156.65'
$ws.Range('E51').Value = '  +0.47%  '
